$p = $ppt.ActivePresentation

# --- Slide 10: "Picks an arm for each learner, accordin to constraints." ---
# -> "Picks a sample for each arm for each learner."
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(2)
$para10 = $sh10.TextFrame.TextRange.Paragraphs(1, 1)

# Remove the trailing ", accordin to constraints" (keep the final period).
$tail10 = $para10.Characters(30, 25)
$tail10.Text = ""

# Expand " an " into " a sample for each " (adds "each" before the existing "arm").
$mid10 = $para10.Characters(6, 4)
$mid10.Text = " a sample for each "

# Re-assert the newly inserted "each" as its own run (forces a run split).
$each10 = $para10.Characters(20, 4)
$each10.Text = "each"

# --- Slide 17: "Alpha values  3,5,10, to be powered by two." ---
# -> "Alpha values  3(2),5(2),10(2)." with superscript "2" exponents.
$s17 = $p.Slides.Item(17)
$sh17 = $s17.Shapes.Item(5)
$para17 = $sh17.TextFrame.TextRange.Paragraphs(5, 1)

# Replace ", to be powered by two" tail (after "3,5,10") with placeholders for the exponents.
$tail17 = $para17.Characters(13, 31)
$tail17.Text = "  3_,5_,10_."

# Turn each placeholder "_" into a superscript "2".
$exp1 = $para17.Characters(16, 1)
$exp1.Text = "2"
$exp1.Font.BaselineOffset = 0.3

$exp2 = $para17.Characters(19, 1)
$exp2.Text = "2"
$exp2.Font.BaselineOffset = 0.3

$exp3 = $para17.Characters(23, 1)
$exp3.Text = "2"
$exp3.Font.BaselineOffset = 0.3
